$wb = $excel.ActiveWorkbook

# Update the status text from "Ready for handoff" to "In Translation"
# on the Overview sheet (columns E/F) and zh-cn/de-de sheets (column C).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"

# Adjust column widths to match the new (narrower) text width.
# NOTE: COM ColumnWidth is character-unit and gets snapped by the host to a
# coarse display-pixel grid, so we feed it the character width whose
# snapped/stored result lands closest to the target stored width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
